$wb = $excel.ActiveWorkbook

# Duplicate the last day's sheet ("03.03") to create the new day's sheet,
# placing the copy right after it -- this carries over the shared layout,
# styles, column widths, frozen panes and formulas for the new "04.03" tab.
$src = $wb.Worksheets.Item("03.03")
$src.Copy([System.Reflection.Missing]::Value, $src)
$new = $wb.Worksheets.Item($src.Index + 1)
$new.Name = "04.03"

# --- Regione data for 04.03 -------------------------------------------------
# Columns: B=Ricoverati con sintomi, C=Terapia Intensiva, D=Isolamento Domiciliare,
#          F=Dimessi Guariti, G=Deceduti, I=Tamponi
# (E=Totale, H=Totale Positivi are SUM formulas and recompute automatically)

# Lombardia (row 2)
$new.Range("B2").Value = 877
$new.Range("C2").Value = 209
$new.Range("D2").Value = 411
$new.Range("F2").Value = 250
$new.Range("G2").Value = 73
$new.Range("I2").Value = 12138

# Veneto (row 3)
$new.Range("B3").Value = 76
$new.Range("C3").Value = 23
$new.Range("D3").Value = 246
$new.Range("F3").Value = 9
$new.Range("G3").Value = 6
$new.Range("I3").Value = 10515

# Emila Romanga (row 4)
$new.Range("B4").Value = 256
$new.Range("C4").Value = 26
$new.Range("D4").Value = 234
$new.Range("F4").Value = 6
$new.Range("G4").Value = 22
$new.Range("I4").Value = 2500

# Piemonte (row 5)
$new.Range("B5").Value = 26
$new.Range("C5").Value = 13
$new.Range("D5").Value = 43
$new.Range("I5").Value = 543

# Liguria (row 6)
$new.Range("B6").Value = 10
$new.Range("C6").Value = 3
$new.Range("D6").Value = 8
$new.Range("I6").Value = 133

# Marche (row 7)
$new.Range("B7").Value = 34
$new.Range("C7").Value = 15
$new.Range("D7").Value = 31
$new.Range("G7").Value = 4
$new.Range("I7").Value = 288

# Toscana (row 8)
$new.Range("B8").Value = 15
$new.Range("C8").Value = 2
$new.Range("D8").Value = 20
$new.Range("I8").Value = 776

# Sicilia (row 9)
$new.Range("B9").Value = 5
$new.Range("D9").Value = 11
$new.Range("I9").Value = 367

# Lazio (row 10)
$new.Range("B10").Value = 15
$new.Range("C10").Value = 3
$new.Range("D10").Value = 9
$new.Range("I10").Value = 995

# Campania (row 11)
$new.Range("D11").Value = 20
$new.Range("I11").Value = 429

# Puglia (row 12)
$new.Range("B12").Value = 4
$new.Range("D12").Value = 3
$new.Range("F12").Value = 1
$new.Range("G12").Value = 1
$new.Range("I12").Value = 322

# Bolzano (row 13) -- unchanged from 03.03

# Calabria (row 14)
$new.Range("B14").Value = 7
$new.Range("D14").ClearContents()
$new.Range("I14").Value = 85

# Sardegna (row 15)
$new.Range("I15").Value = 46

# Umbria (row 16)
$new.Range("D16").Value = 1
$new.Range("I16").Value = 42

# Valle D'Aosta (row 17)
$new.Range("D17").Value = 7
$new.Range("I17").Value = 58

# Friuli V.G. (row 18)
$new.Range("I18").Value = 15

# Trento (row 19)
$new.Range("B19").Value = 3
$new.Range("D19").Value = 15
$new.Range("I19").Value = 376

# Molise (row 20)
$new.Range("D20").Value = 4

# Basilicata (row 21) -- B/C/D unchanged from 03.03
$new.Range("I21").Value = 19

# (row 22) -- B/C/D unchanged from 03.03
$new.Range("I22").Value = 48

# Row 23 totals (B23:I23) are SUM formulas copied from 03.03 -- they
# recompute automatically against the new figures above.

# Make the new sheet the active one, matching Excel's behaviour right after
# a sheet-copy operation, and leave the same cell selected as in the source.
$new.Activate()
$new.Range("I23").Select()
